$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "86×23=" "59×62="
Replace-Text "20×29=" "71×59="
Replace-Text "96×59=" "40×18="
Replace-Text "89×22=" "80×36="
Replace-Text "85×50=" "58×21="
Replace-Text "87×90=" "65×64="
Replace-Text "17×30=" "96×64="
Replace-Text "99×58=" "19×21="
Replace-Text "15×80=" "49×79="
Replace-Text "98×47=" "70×78="
Replace-Text "93×99=" "60×80="
Replace-Text "97×32=" "78×14="
Replace-Text "15×95=" "49×27="
Replace-Text "26×28=" "40×57="
Replace-Text "53×58=" "59×89="
Replace-Text "14×37=" "74×59="
Replace-Text "59×86=" "81×26="
Replace-Text "26×70=" "85×84="
Replace-Text "96×22=" "54×65="
Replace-Text "64×98=" "84×76="
Replace-Text "70×47=" "41×59="
Replace-Text "99×41=" "80×93="
Replace-Text "66×90=" "96×44="
Replace-Text "32×59=" "26×54="
Replace-Text "70×79=" "67×70="
